$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 613
$ws.Range("F4").Value = 6468
$ws.Range("F7").Value = 83
$ws.Range("F10").Value = 27
$ws.Range("F11").Value = 725
$ws.Range("F12").Value = 1206
$ws.Range("F15").Value = 202
$ws.Range("F16").Value = 452
$ws.Range("F18").Value = 24
$ws.Range("F19").Value = 1421
$ws.Range("F20").Value = 678
$ws.Range("F21").Value = 395
$ws.Range("F22").Value = 405
$ws.Range("F23").Value = 84
$ws.Range("F24").Value = 1078
$ws.Range("F25").Value = 170
$ws.Range("F26").Value = 2234
$ws.Range("F27").Value = 261
$ws.Range("F28").Value = 112
$ws.Range("F31").Value = 3632

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 714
$ws.Range("F12").Value = 1023
$ws.Range("F14").Value = 115
$ws.Range("F27").Value = 199
$ws.Range("F31").Value = 214
$ws.Range("F35").Value = 1673

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F7").Value = 435
$ws.Range("F10").Value = 806

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 435
$ws.Range("F7").Value = 806
$ws.Range("F8").Value = 613
$ws.Range("F9").Value = 6468
$ws.Range("F13").Value = 714
$ws.Range("F14").Value = 83
$ws.Range("F17").Value = 27
$ws.Range("F18").Value = 725
$ws.Range("F20").Value = 115
$ws.Range("F21").Value = 115
$ws.Range("F24").Value = 1206
$ws.Range("F27").Value = 202
$ws.Range("F31").Value = 24
$ws.Range("F34").Value = 678
$ws.Range("F35").Value = 395
$ws.Range("F36").Value = 405
$ws.Range("F39").Value = 1078
$ws.Range("F40").Value = 170
$ws.Range("F41").Value = 2234
$ws.Range("F43").Value = 1673
$ws.Range("F44").Value = 1673
$ws.Range("F45").Value = 112
$ws.Range("F47").Value = 3632
